$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = '60.794.09'
$ws.Range("E2").Value = '  -1.27%  '
$ws.Range("D3").Value = '3.386.12'
$ws.Range("E3").Value = '  -1.91%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '''569.38'
$ws.Range("E5").Value = '  -1.63%  '
$ws.Range("D6").Value = '''141.09'
$ws.Range("E6").Value = '  -2.79%  '
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("D8").Value = '3.386.54'
$ws.Range("E8").Value = '  -1.97%  '
$ws.Range("E9").Value = '  -0.49%  '
$ws.Range("D10").Value = '''7.49'
$ws.Range("E10").Value = '  -1.65%  '
$ws.Range("E11").Value = '  -1.57%  '
$ws.Range("D12").Value = '''0.395'
$ws.Range("E12").Value = '  +1.88%  '
$ws.Range("D13").Value = '3.964.55'
$ws.Range("E13").Value = '  -1.87%  '
$ws.Range("D14").Value = '''28.44'
$ws.Range("E14").Value = '  +1.32%  '
$ws.Range("E15").Value = '  +2.14%  '
$ws.Range("E16").Value = '  -1.50%  '
$ws.Range("D17").Value = '3.382.21'
$ws.Range("E17").Value = '  -1.81%  '
$ws.Range("D18").Value = '60.890.60'
$ws.Range("E18").Value = '  -1.35%  '
$ws.Range("D19").Value = '''6.21'
$ws.Range("E19").Value = '  -1.01%  '
$ws.Range("D20").Value = '''13.98'
$ws.Range("E20").Value = '  -1.94%  '
$ws.Range("D21").Value = '''8.98'
$ws.Range("E21").Value = '  -5.83%  '
$ws.Range("D22").Value = '''383.25'
$ws.Range("E22").Value = '  -1.94%  '
$ws.Range("D23").Value = '''0.559'
$ws.Range("E23").Value = '  -1.04%  '
$ws.Range("D24").Value = '''73.73'
$ws.Range("E24").Value = '  +0.58%  '
$ws.Range("E25").Value = '  +0.04%  '
$ws.Range("E26").Value = '  -5.46%  '
$ws.Range("D27").Value = '3.519.59'
$ws.Range("E27").Value = '  -1.95%  '
$ws.Range("E28").Value = '  +0.19%  '
$ws.Range("E29").Value = '  -0.16%  '
$ws.Range("D30").Value = '''7.40'
$ws.Range("E30").Value = '  -3.08%  '
$ws.Range("B32").Value = 'Fetch.AI'
$ws.Range("C32").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D32").Value = '''1.44'
$ws.Range("E32").Value = '  -2.53%  '
$ws.Range("B33").Value = 'PancakeSwap'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D33").Value = '''2.14'
$ws.Range("E33").Value = '  -1.83%  '
$ws.Range("E34").Value = '  +0.03%  '
$ws.Range("D35").Value = '''23.55'
$ws.Range("E35").Value = '  -2.16%  '
$ws.Range("D36").Value = '''6.97'
$ws.Range("E36").Value = '  -0.30%  '
$ws.Range("D37").Value = '''166.18'
$ws.Range("E37").Value = '  -0.49%  '
$ws.Range("D38").Value = '3.416.67'
$ws.Range("E38").Value = '  -1.79%  '
$ws.Range("E39").Value = '  -3.04%  '
$ws.Range("E40").Value = '  -4.60%  '
$ws.Range("D41").Value = '''28.05'
$ws.Range("E41").Value = '  -0.17%  '
$ws.Range("D42").Value = '''0.0770'
$ws.Range("E42").Value = '  -1.28%  '
$ws.Range("D43").Value = '''0.999'
$ws.Range("E43").Value = '  -0.11%  '
$ws.Range("D44").Value = '''0.778'
$ws.Range("E44").Value = '  -3.03%  '
$ws.Range("D45").Value = '''41.92'
$ws.Range("E45").Value = '  -0.87%  '
$ws.Range("D46").Value = '''4.41'
$ws.Range("E46").Value = '  -1.89%  '
$ws.Range("D47").Value = '''1.66'
$ws.Range("E47").Value = '  -4.19%  '
$ws.Range("E48").Value = '  -2.62%  '
$ws.Range("D49").Value = '2.494.36'
$ws.Range("E49").Value = '  -3.79%  '
$ws.Range("D50").Value = '''23.53'
$ws.Range("E50").Value = '  +2.32%  '
$ws.Range("E51").Value = '  -1.65%  '
